$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# fecha1 (row 2) and fecha2 (row 3) sample dates updated for the new climber
$ws.Range("C2").Value = "31_10_2023"
$ws.Range("C3").Value = "02_11_2023"

# DNI numeric value updated for the new climber
$ws.Range("C4").Value = 37504394

# Objetivo text expanded
$ws.Range("C7").Value = "Escalar 7b en la próxima temporada sin lesiones ni molestias y ser feliz"

# Update scroll/selection so the sheet opens scrolled to column C with C8 selected
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C8").Select()
